$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for handback" — both language sheets (zh-cn, de-de) move
# from "awaiting handoff" to "handed back": the status text changes, the
# handoff reason flips from Ignored -> Include, a fresh handback timestamp is
# recorded, and the newly-populated "Latest Target File" / "Latest Handback
# File" columns (E, F) get the file names (with hyperlinks) that were handed
# back.
# ---------------------------------------------------------------------------

function Update-LangSheet($SheetName, $XlfTarget, $HandbackTime, $MdUrlRoot, $XlfUrlRoot) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (B): "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # --- Row 2 (7bf0d961...) ---------------------------------------------
    $ws.Range("E2").Value = "7bf0d961-3b15-40fe-9af6-b20c2f2ff39d.md"
    $ws.Hyperlinks.Add($ws.Range("E2"), "$MdUrlRoot/e2e/7bf0d961-3b15-40fe-9af6-b20c2f2ff39d.md", "", "", "7bf0d961-3b15-40fe-9af6-b20c2f2ff39d.md") | Out-Null
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 15570276

    $ws.Range("F2").Value = "7bf0d961-3b15-40fe-9af6-b20c2f2ff39d.43ebd60e623e5ab3179fd1d280c23504f5c2e864.$XlfTarget"
    $ws.Hyperlinks.Add($ws.Range("F2"), "$XlfUrlRoot/7bf0d961-3b15-40fe-9af6-b20c2f2ff39d.43ebd60e623e5ab3179fd1d280c23504f5c2e864.$XlfTarget", "", "", "7bf0d961-3b15-40fe-9af6-b20c2f2ff39d.43ebd60e623e5ab3179fd1d280c23504f5c2e864.$XlfTarget") | Out-Null
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276

    $ws.Range("G2").Value = $HandbackTime
    $ws.Range("H2").Value = "Include"

    # --- Row 3 (889e6598...) ---------------------------------------------
    $ws.Range("E3").Value = "889e6598-494c-40f9-8f74-51340d8ae09e.md"
    $ws.Hyperlinks.Add($ws.Range("E3"), "$MdUrlRoot/e2e/889e6598-494c-40f9-8f74-51340d8ae09e.md", "", "", "889e6598-494c-40f9-8f74-51340d8ae09e.md") | Out-Null
    $ws.Range("E3").Font.Underline = $true
    $ws.Range("E3").Font.Color = 15570276

    $ws.Range("F3").Value = "889e6598-494c-40f9-8f74-51340d8ae09e.2f9cf2fbda8c6568c97fd13ee3b3cba5a8fde790.$XlfTarget"
    $ws.Hyperlinks.Add($ws.Range("F3"), "$XlfUrlRoot/889e6598-494c-40f9-8f74-51340d8ae09e.2f9cf2fbda8c6568c97fd13ee3b3cba5a8fde790.$XlfTarget", "", "", "889e6598-494c-40f9-8f74-51340d8ae09e.2f9cf2fbda8c6568c97fd13ee3b3cba5a8fde790.$XlfTarget") | Out-Null
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = 15570276

    $ws.Range("G3").Value = $HandbackTime
    $ws.Range("H3").Value = "Include"
}

Update-LangSheet "zh-cn" "zh-cn.xlf" "2016-01-18 05:43:11" `
    "https://github.com/OpenLocalizationTest/oltest/blob/0ae5bff59e02d242c2053d502b06dbf63ee88f18" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/617e963fcda621dc2534ed605509202442611662/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"

Update-LangSheet "de-de" "de-de.xlf" "2016-01-18 05:43:33" `
    "https://github.com/OpenLocalizationTest/oltest/blob/0ae5bff59e02d242c2053d502b06dbf63ee88f18" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a78b03b512647a684309ebff3672fe26db3c1e0d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

Write-Host "Handback report generated"
